# Generate Report for Handoff
# Adds two new "dependency" rows (png files) to each sheet and refreshes
# the existing markdown row with a new handoff guid / timestamps.

$wb = $excel.ActiveWorkbook

$newMd  = "6a0fff05-839d-43f4-8234-a15350d75fef.md"

$png1 = "9b1d5aaa-9ff8-4553-8ea9-1b6e0b3c5328.png"
$png2 = "f8c19475-090e-446a-bd1f-99d3b6b763b3.png"

$readyStatus = "Ready for handoff"
$overviewDate = "2016-49-12 18:49:48"

$zhXlf = "6a0fff05-839d-43f4-8234-a15350d75fef.3dc0f0bd8bc0e64204c7cbb89787e80669d5a39b.zh-cn.xlf"
$deXlf = "6a0fff05-839d-43f4-8234-a15350d75fef.3dc0f0bd8bc0e64204c7cbb89787e80669d5a39b.de-de.xlf"

$zhDate = "2016-03-12 18:49:45"
$deDate = "2016-03-12 18:49:48"

$zeroDate = "0001-01-01 00:00:00"
$includeReason = "Include"
$dependencyReason = "IsDependency"
$dependencyFrom = "e2e\6a0fff05-839d-43f4-8234-a15350d75fef.md"

$png1Target = "5918db4f33c20c0939e7687b9dcb10f2997bca31.png"
$png2Target = "2e7758946861e570accf5d6acc9875ea11fe024c.png"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/5f857608fa9cf1dd12b6f551e30f6e8821970478/e2e/$newMd"
$png1Url = "https://github.com/OpenLocalizationTest/oltest/blob/5f857608fa9cf1dd12b6f551e30f6e8821970478/e2e/$png1"
$png2Url = "https://github.com/OpenLocalizationTest/oltest/blob/5f857608fa9cf1dd12b6f551e30f6e8821970478/e2e/$png2"

$zhXlfUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/215c74d5b98e1c4fa5794f33c0a7dea09446ab81/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$zhPng1Url   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/215c74d5b98e1c4fa5794f33c0a7dea09446ab81/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1Target"
$zhPng2Url   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/215c74d5b98e1c4fa5794f33c0a7dea09446ab81/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2Target"

$deXlfUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40d15eee092e8b0a0a8daafc51e47e1179feefeb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"
$dePng1Url   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40d15eee092e8b0a0a8daafc51e47e1179feefeb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1Target"
$dePng2Url   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40d15eee092e8b0a0a8daafc51e47e1179feefeb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2Target"

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value2 = $newMd
$wsOverview.Range("B2").Value2 = $readyStatus
$wsOverview.Range("C2").Value2 = $readyStatus
$wsOverview.Range("D2").Value2 = $overviewDate

$wsOverview.Range("A3").Value2 = $png1
$wsOverview.Range("B3").Value2 = $readyStatus
$wsOverview.Range("C3").Value2 = $readyStatus
$wsOverview.Range("D3").Value2 = $overviewDate

$wsOverview.Range("A4").Value2 = $png2
$wsOverview.Range("B4").Value2 = $readyStatus
$wsOverview.Range("C4").Value2 = $readyStatus
$wsOverview.Range("D4").Value2 = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $png1Url, "", "", $png1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $png2Url, "", "", $png2) | Out-Null

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

$wsZh.Range("A2").Value2 = $newMd
$wsZh.Range("B2").Value2 = ".md"
$wsZh.Range("C2").Value2 = $readyStatus
$wsZh.Range("D2").Value2 = $zhXlf
$wsZh.Range("E2").Value2 = $zhDate
$wsZh.Range("H2").Value2 = $zeroDate
$wsZh.Range("I2").Value2 = $includeReason

$wsZh.Range("A3").Value2 = $png1
$wsZh.Range("B3").Value2 = ".png"
$wsZh.Range("C3").Value2 = $readyStatus
$wsZh.Range("D3").Value2 = $png1Target
$wsZh.Range("E3").Value2 = $zhDate
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value2 = $zeroDate
$wsZh.Range("I3").Value2 = $dependencyReason
$wsZh.Range("J3").Value2 = $dependencyFrom

$wsZh.Range("A4").Value2 = $png2
$wsZh.Range("B4").Value2 = ".png"
$wsZh.Range("C4").Value2 = $readyStatus
$wsZh.Range("D4").Value2 = $png2Target
$wsZh.Range("E4").Value2 = $zhDate
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value2 = $zeroDate
$wsZh.Range("I4").Value2 = $dependencyReason
$wsZh.Range("J4").Value2 = $dependencyFrom

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, "", "", $zhXlf) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $png1Url, "", "", $png1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $png1Url, "", "", ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhPng1Url, "", "", $png1Target) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $png2Url, "", "", $png2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $png2Url, "", "", ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $zhPng2Url, "", "", $png2Target) | Out-Null

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

$wsDe.Range("A2").Value2 = $newMd
$wsDe.Range("B2").Value2 = ".md"
$wsDe.Range("C2").Value2 = $readyStatus
$wsDe.Range("D2").Value2 = $deXlf
$wsDe.Range("E2").Value2 = $deDate
$wsDe.Range("H2").Value2 = $zeroDate
$wsDe.Range("I2").Value2 = $includeReason

$wsDe.Range("A3").Value2 = $png1
$wsDe.Range("B3").Value2 = ".png"
$wsDe.Range("C3").Value2 = $readyStatus
$wsDe.Range("D3").Value2 = $png1Target
$wsDe.Range("E3").Value2 = $deDate
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value2 = $zeroDate
$wsDe.Range("I3").Value2 = $dependencyReason
$wsDe.Range("J3").Value2 = $dependencyFrom

$wsDe.Range("A4").Value2 = $png2
$wsDe.Range("B4").Value2 = ".png"
$wsDe.Range("C4").Value2 = $readyStatus
$wsDe.Range("D4").Value2 = $png2Target
$wsDe.Range("E4").Value2 = $deDate
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value2 = $zeroDate
$wsDe.Range("I4").Value2 = $dependencyReason
$wsDe.Range("J4").Value2 = $dependencyFrom

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, "", "", $deXlf) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $png1Url, "", "", $png1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $png1Url, "", "", ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $dePng1Url, "", "", $png1Target) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $png2Url, "", "", $png2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $png2Url, "", "", ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $dePng2Url, "", "", $png2Target) | Out-Null

Write-Host "Report generated for handoff"
